$d = $word.ActiveDocument

# Commit message: "Add kee-next to style "dt""
# -> turn on Keep-With-Next for the paragraph style "dt" (this writes
#    <w:keepNext/> into that style's <w:pPr> in word/styles.xml).
$dtStyle = $d.Styles.Item("dt")
$dtStyle.ParagraphFormat.KeepWithNext = $true

# Side-effect observed in the target file: Word's hidden "_GoBack"
# bookmark (tracks the caret position at last save) ends up inside the
# "dt: definition term" paragraph -- right after "dt: definiti" -- since
# that is where the style edit above was made, splitting the run in two.
$r = $d.Content
$found = $r.Find.Execute("dt: definition term", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitAt = $r.Start + 12
    $bmRange = $d.Range($splitAt, $splitAt)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
